$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.663404822349548
$ws.Range("B1").Value = 4.022072792053223
$ws.Range("C1").Value = 3.408957004547119
$ws.Range("D1").Value = 1.550602078437805
$ws.Range("E1").Value = 0.8689704537391663
